$d = $word.ActiveDocument

# Replace "generic" with "diverse" in the spawn-locations sentence.
$d.Content.Find.Execute("generic locations of the players and mazes", $true, $false, $false, $false, $false,
                         $true, 1, $false, "diverse locations of the players and mazes", 2)

# The author's commit split that sentence into three runs around the
# replaced word. Re-select just the new word and nudge a character
# property on/off so the run gets split out from its neighbours.
$r = $d.Content
$r.Find.Execute("diverse", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Bold = 1
$r.Bold = 0
